$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.834.66"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "1.875.16"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.19"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4596"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3876"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07874"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9835"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.71"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "1.855.84"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.986"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.647"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06960"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.11"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009966"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.93"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "28.863.73"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.243"
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("E23").Value = "  -0.55%  "
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.77"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.30"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.992"
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.924"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.25"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09338"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9022"
$ws.Range("E31").Value = "  -2.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.260"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.314"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.266"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.185"
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05759"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02068"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.649"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5647"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1764"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.626"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.247"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.88"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5344"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07036"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.842"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.09"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.511"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.060"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.56"
$ws.Range("E51").Value = "  -0.46%  "
